$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 16-26 -------------------------------------------------
# Values are entered in an order chosen to reproduce the shared-string
# table ordering of the target workbook as closely as possible.

$ws.Range("A16").Value = "_FC_"
$ws.Range("A17").Value = "_BC_"
$ws.Range("A18").Value = "_FSL_"
$ws.Range("C18").Value = "斜正面"
$ws.Range("A19").Value = "_Fin_"
$ws.Range("C19").Value = "內部照"
$ws.Range("A20").Value = "_H1_"
$ws.Range("C20").Value = "細節"
$ws.Range("A21").Value = "_H2_"
$ws.Range("B18").Value = "包包"
$ws.Range("A22").Value = "_Fin_Torso"
$ws.Range("A23").Value = "_Fin_eCom"
$ws.Range("B23").Value = "三合一外套"
$ws.Range("B22").Value = "雙面外套"
$ws.Range("C16").Value = "正面,平拍正面"
$ws.Range("C17").Value = "背面,平拍背面"

# --- C5: retype so the shared string is dropped and recreated -------
$ws.Range("C5").Value = "_TMP_"
$ws.Range("C5").Value = "直立正面"
$ws.Range("C5").Font.Name = "微軟正黑體"

$ws.Range("A24").Value = "_B_Model_"
$ws.Range("C24").Value = "模特背面"
$ws.Range("A25").Value = "_F_Model_"
$ws.Range("C25").Value = "模特正面"
$ws.Range("A26").Value = "_W_Model_"
$ws.Range("C26").Value = "模特側面"

# C21 duplicates the C20 text.
$ws.Range("C21").Value = "細節"

# --- Fonts -----------------------------------------------------------
# A21/A22/A25/A26 use a new plain-Arial font (created before the
# 微軟正黑體 one above so that the cellXfs ordering matches).
$ws.Range("A21").Font.Name = "Arial"

# The remaining "description" cells (column C plus a couple of column B
# cells) reuse the existing Microsoft JhengHei font/style already present
# in the workbook.
$descCells = @("C16","C17","B18","C18","C19","C20","C21","B22","B23","C24","C25","C26")
foreach ($addr in $descCells) {
    $ws.Range($addr).Font.Name = "Microsoft JhengHei"
}

$ws.Range("A22").Font.Name = "Arial"
$ws.Range("A25").Font.Name = "Arial"
$ws.Range("A26").Font.Name = "Arial"

# --- Selection ---------------------------------------------------------
$ws.Range("A5").Select()
